$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Logs")

# Row 1076
$ws.Range("A3:E3").Copy()
$ws.Range("A1076:E1076").PasteSpecial(-4122)
$ws.Range("F936").Copy()
$ws.Range("F1076").PasteSpecial(-4122)
$ws.Range("A1076").Value = 'Pickup Mic'
$ws.Range("B1076").Value = 42704
$ws.Range("C1076").Value = '1600'
$ws.Range("D1076").Value = 'YL'
$ws.Range("E1076").Value = '280M'
$ws.Range("F1076").Value = 'Return IR mic to KT 516 and place battery in charger. We have a key for this room in YKLN 203C'
$ws.Rows.Item(1076).RowHeight = 30

# Row 1077
$ws.Range("A936:E936").Copy()
$ws.Range("A1077:E1077").PasteSpecial(-4122)
$ws.Range("F936").Copy()
$ws.Range("F1077").PasteSpecial(-4122)
$ws.Range("A1077").Value = 'Pickup PC'
$ws.Range("B1077").Value = 42704
$ws.Range("C1077").Value = '1715'
$ws.Range("D1077").Value = 'ATK'
$ws.Range("E1077").Value = '005'
$ws.Range("F1077").Value = 'Return to DB 0003 and plug in for updates !!'

# Row 1078
$ws.Range("A936:E936").Copy()
$ws.Range("A1078:E1078").PasteSpecial(-4122)
$ws.Range("F936").Copy()
$ws.Range("F1078").PasteSpecial(-4122)
$ws.Range("A1078").Value = 'Pickup Projector'
$ws.Range("B1078").Value = 42704
$ws.Range("C1078").Value = '1715'
$ws.Range("D1078").Value = 'ATK'
$ws.Range("E1078").Value = '005'
$ws.Range("F1078").Value = 'Return to ATK 003C'

# Row 1079
$ws.Range("A935:E935").Copy()
$ws.Range("A1079:E1079").PasteSpecial(-4122)
$ws.Range("F7").Copy()
$ws.Range("F1079").PasteSpecial(-4122)
$ws.Range("A1079").Value = 'Setup PC'
$ws.Range("B1079").Value = 42704
$ws.Range("C1079").Value = '1630'
$ws.Range("D1079").Value = 'HNE'
$ws.Range("E1079").Value = '105'
$ws.Range("F1079").Value = 'Equipment from HNES 003'

# Row 1080
$ws.Range("A935:E935").Copy()
$ws.Range("A1080:E1080").PasteSpecial(-4122)
$ws.Range("F7").Copy()
$ws.Range("F1080").PasteSpecial(-4122)
$ws.Range("A1080").Value = 'Setup Projector'
$ws.Range("B1080").Value = 42704
$ws.Range("C1080").Value = '1630'
$ws.Range("D1080").Value = 'HNE'
$ws.Range("E1080").Value = '105'
$ws.Range("F1080").Value = 'Equipment from HNES 003'

# Row 1081
$ws.Range("A935:E935").Copy()
$ws.Range("A1081:E1081").PasteSpecial(-4122)
$ws.Range("F7").Copy()
$ws.Range("F1081").PasteSpecial(-4122)
$ws.Range("A1081").Value = 'Pickup PC'
$ws.Range("B1081").Value = 42704
$ws.Range("C1081").Value = '1730'
$ws.Range("D1081").Value = 'HNE'
$ws.Range("E1081").Value = '105'
$ws.Range("F1081").Value = 'Return to HNES 003'

# Row 1082
$ws.Range("A935:E935").Copy()
$ws.Range("A1082:E1082").PasteSpecial(-4122)
$ws.Range("F936").Copy()
$ws.Range("F1082").PasteSpecial(-4122)
$ws.Range("A1082").Value = 'Pickup Projector'
$ws.Range("B1082").Value = 42704
$ws.Range("C1082").Value = '1730'
$ws.Range("D1082").Value = 'HNE'
$ws.Range("E1082").Value = '105'
$ws.Range("F1082").Value = 'Return to HNES 003'

# Row 1083
$ws.Range("A935:E935").Copy()
$ws.Range("A1083:E1083").PasteSpecial(-4122)
$ws.Range("A1083").Value = 'Demo'
$ws.Range("B1083").Value = 42704
$ws.Range("C1083").Value = '1900'
$ws.Range("D1083").Value = 'SSB'
$ws.Range("E1083").Value = 'W133'

# Row 1084
$ws.Range("A935:E935").Copy()
$ws.Range("A1084:E1084").PasteSpecial(-4122)
$ws.Range("F936").Copy()
$ws.Range("F1084").PasteSpecial(-4122)
$ws.Range("A1084").Value = 'Setup Skype Kit'
$ws.Range("B1084").Value = 42704
$ws.Range("C1084").Value = '1630'
$ws.Range("D1084").Value = 'OSG'
$ws.Range("E1084").Value = '2010'
$ws.Range("F1084").Value = 'Client  is Darren Thorne     darren.thorne@gmail.com    Equipment fro OSG 1014L'
$ws.Rows.Item(1084).RowHeight = 30

# Row 1085
$ws.Range("A935:E935").Copy()
$ws.Range("A1085:E1085").PasteSpecial(-4122)
$ws.Range("F936").Copy()
$ws.Range("F1085").PasteSpecial(-4122)
$ws.Range("A1085").Value = 'Pickup Skype Kit'
$ws.Range("B1085").Value = 42704
$ws.Range("C1085").Value = '1730'
$ws.Range("D1085").Value = 'OSG'
$ws.Range("E1085").Value = '2010'
$ws.Range("F1085").Value = 'Return Skype kit to OSG 1014L'

# Row 1086
$ws.Range("A3:E3").Copy()
$ws.Range("A1086:E1086").PasteSpecial(-4122)
$ws.Range("F7").Copy()
$ws.Range("F1086").PasteSpecial(-4122)
$ws.Range("A1086").Value = 'Pickup Mic'
$ws.Range("B1086").Value = 42704
$ws.Range("C1086").Value = '1630'
$ws.Range("D1086").Value = 'ACW'
$ws.Range("E1086").Value = '206'
$ws.Range("F1086").Value = 'Pick up neck mic labeled "2" with receiver and cables and return to DB 0003'
$ws.Rows.Item(1086).RowHeight = 30

# Row 1087
$ws.Range("A936:E936").Copy()
$ws.Range("A1087:E1087").PasteSpecial(-4122)
$ws.Range("F936").Copy()
$ws.Range("F1087").PasteSpecial(-4122)
$ws.Range("A1087").Value = 'SCLD Student Event'
$ws.Range("B1087").Value = 42704
$ws.Range("C1087").Value = '1900'
$ws.Range("D1087").Value = 'ACE'
$ws.Range("E1087").Value = '004'
$ws.Range("F1087").Value = '752375'

# Row 1088
$ws.Range("A936:E936").Copy()
$ws.Range("A1088:E1088").PasteSpecial(-4122)
$ws.Range("F936").Copy()
$ws.Range("F1088").PasteSpecial(-4122)
$ws.Range("A1088").Value = 'SCLD Student Logout'
$ws.Range("B1088").Value = 42704
$ws.Range("C1088").Value = '2100'
$ws.Range("D1088").Value = 'ACE'
$ws.Range("E1088").Value = '004'
$ws.Range("F1088").Value = '752375'

# Row 1089
$ws.Range("A936:E936").Copy()
$ws.Range("A1089:E1089").PasteSpecial(-4122)
$ws.Range("F936").Copy()
$ws.Range("F1089").PasteSpecial(-4122)
$ws.Range("A1089").Value = 'SCLD Student Event'
$ws.Range("B1089").Value = 42704
$ws.Range("C1089").Value = '1730'
$ws.Range("D1089").Value = 'ACW'
$ws.Range("E1089").Value = '206'
$ws.Range("F1089").Value = 'INC000000754482'

# Row 1090
$ws.Range("A936:E936").Copy()
$ws.Range("A1090:E1090").PasteSpecial(-4122)
$ws.Range("F936").Copy()
$ws.Range("F1090").PasteSpecial(-4122)
$ws.Range("A1090").Value = 'SCLD Student Logout'
$ws.Range("B1090").Value = 42704
$ws.Range("C1090").Value = '2150'
$ws.Range("D1090").Value = 'ACW'
$ws.Range("E1090").Value = '206'
$ws.Range("F1090").Value = 'INC000000754482'

# Row 1091
$ws.Range("A936:E936").Copy()
$ws.Range("A1091:E1091").PasteSpecial(-4122)
$ws.Range("F936").Copy()
$ws.Range("F1091").PasteSpecial(-4122)
$ws.Range("A1091").Value = 'SCLD Student Event'
$ws.Range("B1091").Value = 42704
$ws.Range("C1091").Value = '1730'
$ws.Range("D1091").Value = 'WC'
$ws.Range("E1091").Value = '118'
$ws.Range("F1091").Value = 'INC000000755138'

# Row 1092
$ws.Range("A936:E936").Copy()
$ws.Range("A1092:E1092").PasteSpecial(-4122)
$ws.Range("F936").Copy()
$ws.Range("F1092").PasteSpecial(-4122)
$ws.Range("A1092").Value = 'SCLD Student Logout'
$ws.Range("B1092").Value = 42704
$ws.Range("C1092").Value = '2050'
$ws.Range("D1092").Value = 'WC'
$ws.Range("E1092").Value = '118'
$ws.Range("F1092").Value = 'INC000000755138'

# Row 1097
$ws.Range("A936:E936").Copy()
$ws.Range("A1097:E1097").PasteSpecial(-4122)
$ws.Range("F936").Copy()
$ws.Range("F1097").PasteSpecial(-4122)
$ws.Range("A1097").Value = 'Setup Mic'
$ws.Range("B1097").Value = 42705
$ws.Range("C1097").Value = '1800'
$ws.Range("D1097").Value = 'DB'
$ws.Range("E1097").Value = '2027'
$ws.Range("F1097").Value = 'Neck mic and small PA from DB 0003'

# Row 1098
$ws.Range("A936:E936").Copy()
$ws.Range("A1098:E1098").PasteSpecial(-4122)
$ws.Range("F936").Copy()
$ws.Range("F1098").PasteSpecial(-4122)
$ws.Range("A1098").Value = 'Pickup Mic'
$ws.Range("B1098").Value = 42705
$ws.Range("C1098").Value = '2100'
$ws.Range("D1098").Value = 'DB'
$ws.Range("E1098").Value = '2027'
$ws.Range("F1098").Value = 'Return neck mic and small PA to DB 0003'

# Row 1099
$ws.Range("A935:E935").Copy()
$ws.Range("A1099:E1099").PasteSpecial(-4122)
$ws.Range("A1099").Value = 'Demo'
$ws.Range("B1099").Value = 42705
$ws.Range("C1099").Value = '1900'
$ws.Range("D1099").Value = 'SSB'
$ws.Range("E1099").Value = 'N108'

# Row 1100
$ws.Range("A951:E951").Copy()
$ws.Range("A1100:E1100").PasteSpecial(-4122)
$ws.Range("F7").Copy()
$ws.Range("F1100").PasteSpecial(-4122)
$ws.Range("A1100").Value = 'SCLD Student Event'
$ws.Range("B1100").Value = 42705
$ws.Range("C1100").Value = '1700'
$ws.Range("D1100").Value = 'WC'
$ws.Range("E1100").Value = '118'
$ws.Range("F1100").Value = '752284'

# Row 1101
$ws.Range("A951:E951").Copy()
$ws.Range("A1101:E1101").PasteSpecial(-4122)
$ws.Range("F7").Copy()
$ws.Range("F1101").PasteSpecial(-4122)
$ws.Range("A1101").Value = 'SCLD Student Logout'
$ws.Range("B1101").Value = 42705
$ws.Range("C1101").Value = '2000'
$ws.Range("D1101").Value = 'WC'
$ws.Range("E1101").Value = '118'
$ws.Range("F1101").Value = '752284'

# Row 1102
$ws.Range("A951:E951").Copy()
$ws.Range("A1102:E1102").PasteSpecial(-4122)
$ws.Range("F7").Copy()
$ws.Range("F1102").PasteSpecial(-4122)
$ws.Range("A1102").Value = 'SCLD Student Event'
$ws.Range("B1102").Value = 42705
$ws.Range("C1102").Value = '1900'
$ws.Range("D1102").Value = 'ACE'
$ws.Range("E1102").Value = '004'
$ws.Range("F1102").Value = '752376'

# Row 1103
$ws.Range("A951:E951").Copy()
$ws.Range("A1103:E1103").PasteSpecial(-4122)
$ws.Range("F7").Copy()
$ws.Range("F1103").PasteSpecial(-4122)
$ws.Range("A1103").Value = 'SCLD Student Logout'
$ws.Range("B1103").Value = 42705
$ws.Range("C1103").Value = '2100'
$ws.Range("D1103").Value = 'ACE'
$ws.Range("E1103").Value = '004'
$ws.Range("F1103").Value = '752376'

# Row 1107
$ws.Range("A3:E3").Copy()
$ws.Range("A1107:E1107").PasteSpecial(-4122)
$ws.Range("F7").Copy()
$ws.Range("F1107").PasteSpecial(-4122)
$ws.Range("A1107").Value = 'Pickup Mic'
$ws.Range("B1107").Value = 42706
$ws.Range("C1107").Value = '1700'
$ws.Range("D1107").Value = 'YL'
$ws.Range("E1107").Value = '280N'
$ws.Range("F1107").Value = 'Return 4 IR mics to KT 516 and place batteries in charger    We do have a key for the room in YKLN 203C'
$ws.Rows.Item(1107).RowHeight = 30

# Row 1111
$ws.Range("A3:E3").Copy()
$ws.Range("A1111:E1111").PasteSpecial(-4122)
$ws.Range("A1111").Value = 'Demo'
$ws.Range("B1111").Value = 42709
$ws.Range("C1111").Value = '1900'
$ws.Range("D1111").Value = 'SSB'
$ws.Range("E1111").Value = 'S124'

# Row 1115
$ws.Range("A3:E3").Copy()
$ws.Range("A1115:E1115").PasteSpecial(-4122)
$ws.Range("F7").Copy()
$ws.Range("F1115").PasteSpecial(-4122)
$ws.Range("A1115").Value = 'Setup Mic'
$ws.Range("B1115").Value = 42710
$ws.Range("C1115").Value = '1700'
$ws.Range("D1115").Value = 'SSB'
$ws.Range("E1115").Value = 'W141'
$ws.Range("F1115").Value = 'Podium mic - there/test; one wired audience mic - from rear booth'
$ws.Rows.Item(1115).RowHeight = 30

# Row 1116
$ws.Range("A935:E935").Copy()
$ws.Range("A1116:E1116").PasteSpecial(-4122)
$ws.Range("F7").Copy()
$ws.Range("F1116").PasteSpecial(-4122)
$ws.Range("A1116").Value = 'Demo'
$ws.Range("B1116").Value = 42710
$ws.Range("C1116").Value = '1700'
$ws.Range("D1116").Value = 'SSB'
$ws.Range("E1116").Value = 'W141'
$ws.Range("F1116").Value = 'Neck mic and PC'

# Row 1117
$ws.Range("A935:E935").Copy()
$ws.Range("A1117:E1117").PasteSpecial(-4122)
$ws.Range("F7").Copy()
$ws.Range("F1117").PasteSpecial(-4122)
$ws.Range("A1117").Value = 'Pickup Mic'
$ws.Range("B1117").Value = 42710
$ws.Range("C1117").Value = '2000'
$ws.Range("D1117").Value = 'SSB'
$ws.Range("E1117").Value = 'W141'
$ws.Range("F1117").Value = 'Pick up oneaudience mic , stand and cable and return to rear booth'
$ws.Rows.Item(1117).RowHeight = 30

# Row 1122
$ws.Range("A3:E3").Copy()
$ws.Range("A1122:E1122").PasteSpecial(-4122)
$ws.Range("F7").Copy()
$ws.Range("F1122").PasteSpecial(-4122)
$ws.Range("A1122").Value = 'Pickup Mic'
$ws.Range("B1122").Value = 42711
$ws.Range("C1122").Value = '1700'
$ws.Range("D1122").Value = 'KT'
$ws.Range("E1122").Value = '519'
$ws.Range("F1122").Value = 'Return mic, stand and cable to KT 516 / we have a key for the room in KT 516'
$ws.Rows.Item(1122).RowHeight = 30

$ws.Application.CutCopyMode = $false

$ws.Range("F1122").Select()
